# Refresh the crypto price/volume table (cols D = Price, E = Volume(1h)).
# Values are stored as text in the source sheet (e.g. "67.538.67" uses dots
# as thousands separators, so it can't be a real number), and some of the
# replacement figures look numeric (e.g. "6.59"). A leading apostrophe
# forces Excel to keep those as literal text instead of coercing them to
# a floating point number, matching the original inlineStr cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.538.67"
$ws.Range("E2").Value = "'  -0.62%  "

$ws.Range("D3").Value = "'3.228.96"
$ws.Range("E3").Value = "'  -0.91%  "

$ws.Range("E4").Value = "'  +0.02%  "

$ws.Range("D5").Value = "'579.87"

$ws.Range("D6").Value = "'182.31"
$ws.Range("E6").Value = "'  -0.91%  "

$ws.Range("E7").Value = "'  +0.02%  "

$ws.Range("E8").Value = "'  +0.70%  "

$ws.Range("D9").Value = "'3.228.50"
$ws.Range("E9").Value = "'  -0.84%  "

$ws.Range("E10").Value = "'  -2.79%  "

$ws.Range("D11").Value = "'6.59"
$ws.Range("E11").Value = "'  -1.37%  "

$ws.Range("E12").Value = "'  -0.70%  "

$ws.Range("D13").Value = "'3.790.48"
$ws.Range("E13").Value = "'  -0.91%  "

$ws.Range("E14").Value = "'  +0.03%  "

$ws.Range("D15").Value = "'27.71"
$ws.Range("E15").Value = "'  -3.00%  "

$ws.Range("D16").Value = "'67.585.45"
$ws.Range("E16").Value = "'  -0.61%  "

$ws.Range("D17").Value = "'0.0000169"
$ws.Range("E17").Value = "'  -1.99%  "

$ws.Range("D18").Value = "'3.230.45"
$ws.Range("E18").Value = "'  -0.77%  "

$ws.Range("D19").Value = "'5.74"
$ws.Range("E19").Value = "'  -1.79%  "

$ws.Range("D20").Value = "'13.43"
$ws.Range("E20").Value = "'  -1.23%  "

$ws.Range("D21").Value = "'395.33"
$ws.Range("E21").Value = "'  +3.67%  "

$ws.Range("D22").Value = "'7.53"
$ws.Range("E22").Value = "'  -1.85%  "

$ws.Range("E23").Value = "'  +0.05%  "

$ws.Range("D24").Value = "'70.85"
$ws.Range("E24").Value = "'  -0.64%  "

$ws.Range("E25").Value = "'  -0.33%  "

$ws.Range("E26").Value = "'  -2.57%  "

$ws.Range("E27").Value = "'  +1.92%  "

$ws.Range("D28").Value = "'9.55"
$ws.Range("E28").Value = "'  -2.99%  "

$ws.Range("E29").Value = "'  +0.17%  "

$ws.Range("E30").Value = "'  -1.62%  "

$ws.Range("D31").Value = "'5.57"
$ws.Range("E31").Value = "'  -1.91%  "

$ws.Range("D32").Value = "'22.62"
$ws.Range("E32").Value = "'  -1.16%  "

$ws.Range("E33").Value = "'  -4.63%  "

$ws.Range("E35").Value = "'  -1.41%  "

$ws.Range("D36").Value = "'161.44"
$ws.Range("E36").Value = "'  -0.96%  "

$ws.Range("E37").Value = "'  -5.25%  "

$ws.Range("E38").Value = "'  +1.47%  "

$ws.Range("D39").Value = "'26.39"
$ws.Range("E39").Value = "'  -0.92%  "

$ws.Range("D40").Value = "'0.804"
$ws.Range("E40").Value = "'  -3.59%  "

$ws.Range("D41").Value = "'4.57"
$ws.Range("E41").Value = "'  -0.77%  "

$ws.Range("D42").Value = "'6.48"
$ws.Range("E42").Value = "'  -4.50%  "

$ws.Range("D43").Value = "'2.46"
$ws.Range("E43").Value = "'  -5.20%  "

$ws.Range("E44").Value = "'  -1.03%  "

$ws.Range("D45").Value = "'40.64"
$ws.Range("E45").Value = "'  -1.85%  "

$ws.Range("D46").Value = "'2.607.16"
$ws.Range("E46").Value = "'  -1.29%  "

$ws.Range("D47").Value = "'24.67"
$ws.Range("E47").Value = "'  -2.91%  "

$ws.Range("D48").Value = "'335.02"
$ws.Range("E48").Value = "'  -3.17%  "

$ws.Range("E49").Value = "'  -2.82%  "

$ws.Range("D50").Value = "'6.27"
$ws.Range("E50").Value = "'  +0.45%  "

$ws.Range("E51").Value = "'  -1.65%  "
